$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 2: FAPs | Ccl12 | Ackr4 | ECs ----
$ws.Cells.Item(2,9).Value  = 0.005723000769734084
$ws.Cells.Item(2,10).Value = 0.005723000769734084
$ws.Cells.Item(2,13).Value = 0.377371
$ws.Cells.Item(2,14).Value = 1.132113
$ws.Cells.Item(2,15).Value = 0.4698794580655765
$ws.Cells.Item(2,16).Value = 0.4698794580655764
$ws.Cells.Item(2,17).Value = 0.1064769887146667
$ws.Cells.Item(2,18).Value = 0.9582928984319999
$ws.Cells.Item(2,19).Value = 0.002689120500191529
$ws.Cells.Item(2,20).Value = 0.002689120500191528

# ---- Row 3: FAPs | Ccl12 | Ackr4 | FAPs ----
$ws.Cells.Item(3,9).Value  = 0.005723000769734084
$ws.Cells.Item(3,10).Value = 0.005723000769734084
$ws.Cells.Item(3,15).Value = 0.443286188209444
$ws.Cells.Item(3,16).Value = 0.443286188209444
$ws.Cells.Item(3,17).Value = 0.1004508233955556
$ws.Cells.Item(3,18).Value = 0.9040574105600001
$ws.Cells.Item(3,19).Value = 0.002536927196335137
$ws.Cells.Item(3,20).Value = 0.002536927196335137

# ---- Row 4: FAPs | Ccl12 | Ackr4 | MuSCs ----
$ws.Cells.Item(4,9).Value  = 0.005723000769734084
$ws.Cells.Item(4,10).Value = 0.005723000769734084
$ws.Cells.Item(4,13).Value = 0.06973866666666667
$ws.Cells.Item(4,14).Value = 0.209216
$ws.Cells.Item(4,15).Value = 0.08683435372497944
$ws.Cells.Item(4,16).Value = 0.08683435372497944
$ws.Cells.Item(4,17).Value = 0.01967709024711111
$ws.Cells.Item(4,18).Value = 0.177093812224
$ws.Cells.Item(4,19).Value = 0.0004969530732074191
$ws.Cells.Item(4,20).Value = 0.0004969530732074191

# ---- Row 5: Inflammatory-Mac | Ccl12 | Ackr4 | ECs (was FAPs|...|Resolving-Mac) ----
$ws.Cells.Item(5,1).Value  = "Inflammatory-Mac"
$ws.Cells.Item(5,4).Value  = "ECs"
$ws.Cells.Item(5,5).Value  = 3
$ws.Cells.Item(5,6).Value  = 1
$ws.Cells.Item(5,7).Value  = 30.199365
$ws.Cells.Item(5,8).Value  = 90.598095
$ws.Cells.Item(5,9).Value  = 0.6125398923302606
$ws.Cells.Item(5,10).Value = 0.6125398923302606
$ws.Cells.Item(5,13).Value = 0.377371
$ws.Cells.Item(5,14).Value = 1.132113
$ws.Cells.Item(5,15).Value = 0.4698794580655765
$ws.Cells.Item(5,16).Value = 0.4698794580655764
$ws.Cells.Item(5,17).Value = 11.396364569415
$ws.Cells.Item(5,18).Value = 102.567281124735
$ws.Cells.Item(5,19).Value = 0.2878199126516894
$ws.Cells.Item(5,20).Value = 0.2878199126516894

# ---- Row 6: Inflammatory-Mac | Ccl12 | Ackr4 | FAPs (was ...|ECs) ----
$ws.Cells.Item(6,4).Value  = "FAPs"
$ws.Cells.Item(6,7).Value  = 30.199365
$ws.Cells.Item(6,8).Value  = 90.598095
$ws.Cells.Item(6,9).Value  = 0.6125398923302606
$ws.Cells.Item(6,10).Value = 0.6125398923302606
$ws.Cells.Item(6,13).Value = 0.3560133333333333
$ws.Cells.Item(6,14).Value = 1.06804
$ws.Cells.Item(6,15).Value = 0.443286188209444
$ws.Cells.Item(6,16).Value = 0.443286188209444
$ws.Cells.Item(6,17).Value = 10.7513765982
$ws.Cells.Item(6,18).Value = 96.76238938380001
$ws.Cells.Item(6,19).Value = 0.2715304739973045
$ws.Cells.Item(6,20).Value = 0.2715304739973045

# ---- Row 7: Inflammatory-Mac | Ccl12 | Ackr4 | MuSCs (was ...|FAPs) ----
$ws.Cells.Item(7,4).Value  = "MuSCs"
$ws.Cells.Item(7,7).Value  = 30.199365
$ws.Cells.Item(7,8).Value  = 90.598095
$ws.Cells.Item(7,9).Value  = 0.6125398923302606
$ws.Cells.Item(7,10).Value = 0.6125398923302606
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.06973866666666667
$ws.Cells.Item(7,14).Value = 0.209216
$ws.Cells.Item(7,15).Value = 0.08683435372497944
$ws.Cells.Item(7,16).Value = 0.08683435372497944
$ws.Cells.Item(7,17).Value = 2.10606344928
$ws.Cells.Item(7,18).Value = 18.95457104352
$ws.Cells.Item(7,19).Value = 0.05318950568126667
$ws.Cells.Item(7,20).Value = 0.05318950568126667

# ---- Row 8: Resolving-Mac | Ccl12 | Ackr4 | ECs (was Inflammatory-Mac|...|MuSCs) ----
$ws.Cells.Item(8,1).Value  = "Resolving-Mac"
$ws.Cells.Item(8,4).Value  = "ECs"
$ws.Cells.Item(8,7).Value  = 18.820355
$ws.Cells.Item(8,8).Value  = 56.461065
$ws.Cells.Item(8,9).Value  = 0.3817371069000054
$ws.Cells.Item(8,10).Value = 0.3817371069000054
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.377371
$ws.Cells.Item(8,14).Value = 1.132113
$ws.Cells.Item(8,15).Value = 0.4698794580655765
$ws.Cells.Item(8,16).Value = 0.4698794580655764
$ws.Cells.Item(8,17).Value = 7.102256186704999
$ws.Cells.Item(8,18).Value = 63.920305680345
$ws.Cells.Item(8,19).Value = 0.1793704249136956
$ws.Cells.Item(8,20).Value = 0.1793704249136955

# ---- Row 9: Resolving-Mac | Ccl12 | Ackr4 | FAPs (was Inflammatory-Mac|...|Resolving-Mac) ----
$ws.Cells.Item(9,1).Value  = "Resolving-Mac"
$ws.Cells.Item(9,4).Value  = "FAPs"
$ws.Cells.Item(9,7).Value  = 18.820355
$ws.Cells.Item(9,8).Value  = 56.461065
$ws.Cells.Item(9,9).Value  = 0.3817371069000054
$ws.Cells.Item(9,10).Value = 0.3817371069000054
$ws.Cells.Item(9,13).Value = 0.3560133333333333
$ws.Cells.Item(9,14).Value = 1.06804
$ws.Cells.Item(9,15).Value = 0.443286188209444
$ws.Cells.Item(9,16).Value = 0.443286188209444
$ws.Cells.Item(9,17).Value = 6.700297318066666
$ws.Cells.Item(9,18).Value = 60.3026758626
$ws.Cells.Item(9,19).Value = 0.1692187870158045
$ws.Cells.Item(9,20).Value = 0.1692187870158045

# ---- Row 10: Resolving-Mac | Ccl12 | Ackr4 | MuSCs (was ...|ECs) ----
$ws.Cells.Item(10,4).Value  = "MuSCs"
$ws.Cells.Item(10,7).Value  = 18.820355
$ws.Cells.Item(10,8).Value  = 56.461065
$ws.Cells.Item(10,9).Value  = 0.3817371069000054
$ws.Cells.Item(10,10).Value = 0.3817371069000054
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.06973866666666667
$ws.Cells.Item(10,14).Value = 0.209216
$ws.Cells.Item(10,15).Value = 0.08683435372497944
$ws.Cells.Item(10,16).Value = 0.08683435372497944
$ws.Cells.Item(10,17).Value = 1.312506463893333
$ws.Cells.Item(10,18).Value = 11.81255817504
$ws.Cells.Item(10,19).Value = 0.03314789497050536
$ws.Cells.Item(10,20).Value = 0.03314789497050536

# ---- Remove the now-obsolete trailing rows 11-13 ----
$ws.Range("A11:T13").Delete() | Out-Null
